$wb = $excel.ActiveWorkbook

# --- Horizontal_Data sheet: reshape the yearly rows from one packed
# "csv-in-a-cell" string per year into individual month cells, and give
# the month-name header row (B1:M1) the same header style as A1/the
# vertical table's header column, so header cells look consistent
# regardless of whether the table is laid out vertically or horizontally.
$ws = $wb.Worksheets.Item("Horizontal_Data")

# Data rows: year label in column A, then the twelve monthly figures.
$data = @{
    2 = @("1958", " 340", " 318", " 362", " 348", " 363", " 435", " 491", " 505", " 404", " 359", " 310", " 337")
    3 = @("1959", " 360", " 342", " 406", " 396", " 420", " 472", " 548", " 559", " 463", " 407", " 362", " 405")
    4 = @("1960", " 417", " 391", " 419", " 461", " 472", " 535", " 622", " 606", " 508", " 461", " 390", " 432")
}

foreach ($r in $data.Keys) {
    $rowRange = $ws.Range($ws.Cells.Item($r, 1), $ws.Cells.Item($r, 13))
    # Force text storage so the leading space / bare year string survives
    # instead of Excel auto-coercing " 340" -> 340.
    $rowRange.NumberFormat = "@"
    $values = $data[$r]
    for ($i = 0; $i -lt $values.Length; $i++) {
        $ws.Cells.Item($r, $i + 1).Value = $values[$i]
    }
    # Drop the temporary text format again so the cells end up unstyled,
    # same as the rest of the data area.
    $rowRange.ClearFormats()
}

# Give the month-name header row the same header style as A1 (s="2").
$ws.Range("A1").Copy()
$ws.Range("B1:M1").PasteSpecial(-4122)

# Column A no longer needs to be wide enough to host the packed string;
# shrink it back down to a normal column width.
$ws.Columns.Item(1).ColumnWidth = 7.5

# The old leftover placeholder cell in A5 goes away now that A2:A4 hold
# real data instead of a single summary cell.
$ws.Cells.Item(5, 1).Clear()
